$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 356.23077
$ws.Range("I4").Value = 304.25
$ws.Range("K4").Value = 304.25
$ws.Range("M4").Value = -190.25

$ws.Range("H5").Value = 384.53845
$ws.Range("I5").Value = 377.66666
$ws.Range("K5").Value = 377.66666
$ws.Range("M5").Value = -262.66666

$ws.Range("H9").Value = 135.5
$ws.Range("I9").Value = 120.8
$ws.Range("K9").Value = 120.8
$ws.Range("M9").Value = 48.2

$ws.Range("H17").Value = 615.2381
$ws.Range("I17").Value = 200
$ws.Range("K17").Value = 600
$ws.Range("M17").Value = -432

$ws.Range("H39").Value = 57.3
$ws.Range("I39").Value = 55.333332
$ws.Range("K39").Value = 165.999996
$ws.Range("M39").Value = 130.000004

$ws.Range("H42").Value = 572.2727
$ws.Range("I42").Value = 411.875
$ws.Range("J42").Value = 1000
$ws.Range("K42").Value = 1235.625
$ws.Range("L42").Value = 3000
$ws.Range("M42").Value = -1005.625
$ws.Range("N42").Value = -3460

$ws.Range("H98").Value = 4419.909
$ws.Range("I98").Value = 4644.1577
$ws.Range("J98").Value = 2999.6667
$ws.Range("K98").Value = 4644.1577
$ws.Range("L98").Value = 2999.6667
$ws.Range("M98").Value = -3146.1577
$ws.Range("N98").Value = -5995.6667

$ws.Range("H122").Value = 4419.909
$ws.Range("I122").Value = 4644.1577
$ws.Range("J122").Value = 2999.6667
$ws.Range("K122").Value = 13932.4731
$ws.Range("L122").Value = 8999.000100000001
$ws.Range("M122").Value = -11482.4731
$ws.Range("N122").Value = -13899.0001

$ws.Range("H127").Value = 1094.9565
$ws.Range("I127").Value = 611.5
$ws.Range("K127").Value = 1834.5
$ws.Range("M127").Value = 3125.5

$ws.Range("H133").Value = 32995.8
$ws.Range("J133").Value = 32995.8
$ws.Range("L133").Value = 32995.8
$ws.Range("N133").Value = -43115.8

$ws.Range("H134").Value = 50663.332
$ws.Range("J134").Value = 50663.332
$ws.Range("L134").Value = 50663.332
$ws.Range("N134").Value = -60803.332

$ws.Range("H137").Value = 1061.6747
$ws.Range("I137").Value = 885.7111
$ws.Range("K137").Value = 2657.1333
$ws.Range("M137").Value = -107.1333

$ws.Range("H141").Value = 524.62
$ws.Range("I141").Value = 524.62
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 1573.86
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 3606.14
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 858.1142599999999
$ws.Range("I74").Value = 741.59375
$ws.Range("K74").Value = 741.59375
$ws.Range("M74").Value = 132.40625

$ws.Range("H77").Value = 858.1142599999999
$ws.Range("I77").Value = 741.59375
$ws.Range("K77").Value = 3707.96875
$ws.Range("M77").Value = 660.03125

$ws.Range("H132").Value = 2083.2708
$ws.Range("I132").Value = 1982.4286
$ws.Range("K132").Value = 5947.2858
$ws.Range("M132").Value = -3417.2858

$ws.Range("H138").Value = 50240
$ws.Range("J138").Value = 50240
$ws.Range("L138").Value = 50240
$ws.Range("N138").Value = -60520

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 10000
$ws.Range("I26").Value = 10000
$ws.Range("K26").Value = 10000
$ws.Range("M26").Value = -9708

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H94").Value = 62501504
$ws.Range("I94").Value = 250000000
$ws.Range("J94").Value = 2003.3334
$ws.Range("K94").Value = 250000000
$ws.Range("L94").Value = 2003.3334
$ws.Range("M94").Value = -249999549
$ws.Range("N94").Value = -2905.3334

$ws.Range("H134").Value = 4235.1523
$ws.Range("I134").Value = 1492.4054
$ws.Range("J134").Value = 15510.889
$ws.Range("K134").Value = 4477.216200000001
$ws.Range("L134").Value = 46532.667
$ws.Range("M134").Value = -1942.216200000001
$ws.Range("N134").Value = -51602.667

$ws.Range("H140").Value = 20780
$ws.Range("J140").Value = 20780
$ws.Range("L140").Value = 20780
$ws.Range("N140").Value = -31140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 58845
$ws.Range("I22").Value = 770
$ws.Range("J22").Value = 87882.5
$ws.Range("K22").Value = 770
$ws.Range("L22").Value = 87882.5
$ws.Range("M22").Value = -420
$ws.Range("N22").Value = -88582.5

$ws.Range("H31").Value = 2460.6667
$ws.Range("I31").Value = 2679.4666
$ws.Range("K31").Value = 2679.4666
$ws.Range("M31").Value = -2384.4666

$ws.Range("H34").Value = 2460.6667
$ws.Range("I34").Value = 2679.4666
$ws.Range("K34").Value = 2679.4666
$ws.Range("M34").Value = -2477.4666

$ws.Range("H105").Value = 1016
$ws.Range("I105").Value = 793.3333
$ws.Range("J105").Value = 1350
$ws.Range("K105").Value = 793.3333
$ws.Range("L105").Value = 1350
$ws.Range("M105").Value = 953.6667
$ws.Range("N105").Value = -4844

$ws.Range("H108").Value = 19154.75
$ws.Range("J108").Value = 18666
$ws.Range("L108").Value = 18666
$ws.Range("N108").Value = -26346

$ws.Range("H132").Value = 3365.182
$ws.Range("I132").Value = 3260.1914
$ws.Range("K132").Value = 9780.574200000001
$ws.Range("M132").Value = -7250.574200000001

$ws.Range("H134").Value = 8334323.5
$ws.Range("I134").Value = 998.2075
$ws.Range("K134").Value = 2994.6225
$ws.Range("M134").Value = -459.6224999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1785.3043
$ws.Range("J5").Value = 1300
$ws.Range("L5").Value = 3900
$ws.Range("N5").Value = -4124

$ws.Range("H121").Value = 608
$ws.Range("J121").Value = 899.3333
$ws.Range("L121").Value = 2697.9999
$ws.Range("N121").Value = -5317.9999

$ws.Range("H133").Value = 2304.2

$ws.Range("H135").Value = 1785.3043
$ws.Range("J135").Value = 1300
$ws.Range("L135").Value = 11700
$ws.Range("N135").Value = -16770

$ws.Range("H137").Value = 22729594
$ws.Range("J137").Value = 3643.7646
$ws.Range("L137").Value = 10931.2938
$ws.Range("N137").Value = -21131.2938

$ws.Range("H139").Value = 1655.3549
$ws.Range("I139").Value = 1639.8823
$ws.Range("J139").Value = 1674.1428
$ws.Range("K139").Value = 4919.6469
$ws.Range("L139").Value = 5022.428400000001
$ws.Range("M139").Value = 220.3531000000003
$ws.Range("N139").Value = -15302.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1802.8077
$ws.Range("J126").Value = 2138.5
$ws.Range("L126").Value = 6415.5
$ws.Range("N126").Value = -11355.5

$ws.Range("H132").Value = 1777.1892
$ws.Range("I132").Value = 1556.1111
$ws.Range("J132").Value = 2374.1
$ws.Range("K132").Value = 4668.3333
$ws.Range("L132").Value = 7122.299999999999
$ws.Range("M132").Value = -2138.3333
$ws.Range("N132").Value = -12182.3

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1923.75
$ws.Range("I7").Value = 1923.75
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1923.75
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1811.75
$ws.Range("N7").ClearContents()

$ws.Range("H22").Value = 1232
$ws.Range("I22").Value = 1122.5
$ws.Range("K22").Value = 1122.5
$ws.Range("M22").Value = -827.5

$ws.Range("H27").Value = 1232
$ws.Range("I27").Value = 1122.5
$ws.Range("K27").Value = 1122.5
$ws.Range("M27").Value = -1015.5

$ws.Range("H122").Value = 41691230
$ws.Range("I122").Value = 125050000
$ws.Range("K122").Value = 375150000
$ws.Range("M122").Value = -375147550

$ws.Range("H126").Value = 1923.75
$ws.Range("I126").Value = 1923.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5771.25
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3301.25
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 18513.9
$ws.Range("I132").Value = 1422.8055
$ws.Range("K132").Value = 4268.416499999999
$ws.Range("M132").Value = -1738.416499999999

$ws.Range("H134").Value = 31666.666
$ws.Range("J134").Value = 31666.666
$ws.Range("L134").Value = 31666.666
$ws.Range("N134").Value = -41806.666

$ws.Range("H136").Value = 3175.4546
$ws.Range("I136").Value = 3183.8333
$ws.Range("K136").Value = 9551.499899999999
$ws.Range("M136").Value = -7001.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 21501500
$ws.Range("J5").Value = 21501500
$ws.Range("L5").Value = 21501500
$ws.Range("N5").Value = -21501724

$ws.Range("H132").Value = 1733.4667
$ws.Range("I132").Value = 1694.5471
$ws.Range("K132").Value = 5083.6413
$ws.Range("M132").Value = -2553.6413

$ws.Range("H133").Value = 45500
$ws.Range("J133").Value = 45500
$ws.Range("L133").Value = 45500
$ws.Range("N133").Value = -55620

$ws.Range("H136").Value = 502.98038
$ws.Range("I136").Value = 394.5
$ws.Range("K136").Value = 1183.5
$ws.Range("M136").Value = 1366.5
